# Added stats2.OptIntVal into PiggyBank
# The PiggyBank stats table (A2:F25) gets re-sorted/rebuilt once the new
# stats2.OptIntVal entry is folded in. All per-key (column A) rows keep
# their original B:F values, only their row position changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(901,16,15,45,60,60),
    @(401,9,48,67,75,45),
    @(601,9,60,67,60,42),
    @(801,3,67,65,52,45),
    @(902,1,0,0,0,0),
    @(1001,18,30,75,60,72),
    @(301,6,45,30,60,45),
    @(1202,2,10,10,10,10),
    @(1201,2,10,10,10,10),
    @(1203,3,15,15,15,15),
    @(101,9,30,15,60,15),
    @(501,9,52,30,75,45),
    @(701,3,90,45,97,15),
    @(201,9,30,15,45,30),
    @(802,0,4,5,4,0),
    @(2,0,2,2,2,2),
    @(1101,0,15,30,30,0),
    @(3,0,3,3,3,3),
    @(502,0,4,0,0,0),
    @(1,0,2,2,2,2),
    @(402,0,0,4,0,0),
    @(602,0,0,4,0,9),
    @(702,0,0,0,4,0),
    @(1002,0,0,0,0,9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
